$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.148.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.274.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.04"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.88"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.12"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.616.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.862"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.274.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.166.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.91"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.34"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.46"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0898"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.26%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.129"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0371"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.89"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.29"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "76.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.37%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.237"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.07"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.51"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0996"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.36"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.598"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.33%  "
